$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header columns (AD/AE/AF), copying the existing header
# style (bold, centered, bordered) from the adjacent "Unnamed: 28" header
# cell (AC1) so the new headers match the formatting of the rest of row 1.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"
$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
$ws.Range("AD2:AD51").Value = 68
$ws.Range("AE2:AE51").Value = 94
$ws.Range("AF2:AF51").Value = 0
